$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Zero out the total_change (G), hitting_change (H) and pitching_change (I)
# columns for every data row (rows 2-15).
$ws.Range("G2:I15").Value = 0.0
